$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed TPM-derived values for rows 2-4 (columns M..T)
$ws.Range("M2").Value = 0.274713
$ws.Range("N2").Value = 0.824139
$ws.Range("O2").Value = 0.1055967877339779
$ws.Range("P2").Value = 0.1055967877339779
$ws.Range("Q2").Value = 0.015119745665
$ws.Range("R2").Value = 0.136077710985
$ws.Range("S2").Value = 0.1055967877339779
$ws.Range("T2").Value = 0.1055967877339779

$ws.Range("O3").Value = 0.3415960415058637
$ws.Range("P3").Value = 0.3415960415058638
$ws.Range("S3").Value = 0.3415960415058637
$ws.Range("T3").Value = 0.3415960415058638

$ws.Range("M4").Value = 1.438143333333333
$ws.Range("N4").Value = 4.31443
$ws.Range("O4").Value = 0.5528071707601584
$ws.Range("P4").Value = 0.5528071707601584
$ws.Range("Q4").Value = 0.07915301216111112
$ws.Range("R4").Value = 0.71237710945
$ws.Range("S4").Value = 0.5528071707601584
$ws.Range("T4").Value = 0.5528071707601584

# Remove the row for the "Resolving-Mac" target cluster entirely
$ws.Rows.Item(5).Delete()
